$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(42, 8).Value = 465.07144  # H42: was 469.42856
$ws.Cells.Item(42, 9).Value = 49.3  # I42: was 53.77778
$ws.Cells.Item(42, 10).Value = 1504.5  # J42: was 1217.6
$ws.Cells.Item(42, 11).Value = 147.9  # K42: was 161.33334
$ws.Cells.Item(42, 12).Value = 4513.5  # L42: was 3652.8
$ws.Cells.Item(42, 13).Value = 82.10000000000002  # M42: was 68.66666000000001
$ws.Cells.Item(42, 14).Value = -4973.5  # N42: was -4112.799999999999
$ws.Cells.Item(69, 8).Value = 3124  # H69: was 3014.9167
$ws.Cells.Item(69, 9).Value = 0  # I69: was 2000
$ws.Cells.Item(69, 10).Value = 3124  # J69: was 3107.182
$ws.Cells.Item(69, 11).Value = 0  # K69: was 6000
$ws.Cells.Item(69, 12).Value = 9372  # L69: was 9321.545999999998
$ws.Cells.Item(69, 13).ClearContents()  # M69: was -5126
$ws.Cells.Item(69, 14).Value = -11120  # N69: was -11069.546
$ws.Cells.Item(72, 8).Value = 3124  # H72: was 3014.9167
$ws.Cells.Item(72, 9).Value = 0  # I72: was 2000
$ws.Cells.Item(72, 10).Value = 3124  # J72: was 3107.182
$ws.Cells.Item(72, 11).Value = 0  # K72: was 18000
$ws.Cells.Item(72, 12).Value = 28116  # L72: was 27964.638
$ws.Cells.Item(72, 13).ClearContents()  # M72: was -13632
$ws.Cells.Item(72, 14).Value = -36852  # N72: was -36700.638
$ws.Cells.Item(112, 8).Value = 1504.174  # H112: was 1483.1666
$ws.Cells.Item(112, 10).Value = 1756  # J112: was 1711.5294
$ws.Cells.Item(112, 12).Value = 5268  # L112: was 5134.5882
$ws.Cells.Item(112, 14).Value = -7484  # N112: was -7350.5882
$ws.Cells.Item(125, 8).Value = 657.0909  # H125: was 500.2
$ws.Cells.Item(125, 9).Value = 458  # I125: was 433.66666
$ws.Cells.Item(125, 10).Value = 770.8570999999999  # J125: was 600
$ws.Cells.Item(125, 11).Value = 4122  # K125: was 3902.99994
$ws.Cells.Item(125, 12).Value = 6937.7139  # L125: was 5400
$ws.Cells.Item(125, 13).Value = -1662  # M125: was -1442.99994
$ws.Cells.Item(125, 14).Value = -11857.7139  # N125: was -10320
$ws.Cells.Item(127, 8).Value = 855.86365  # H127: was 877.5714
$ws.Cells.Item(127, 9).Value = 580.7778  # I127: was 603.375
$ws.Cells.Item(127, 11).Value = 1742.3334  # K127: was 1810.125
$ws.Cells.Item(127, 13).Value = 3217.6666  # M127: was 3149.875
$ws.Cells.Item(129, 8).Value = 821.5714  # H129: was 832.8
$ws.Cells.Item(129, 10).Value = 1034.1  # J129: was 1030.091
$ws.Cells.Item(129, 12).Value = 3102.3  # L129: was 3090.273
$ws.Cells.Item(129, 14).Value = -13102.3  # N129: was -13090.273
$ws.Cells.Item(137, 8).Value = 1425.8431  # H137: was 1528.6666
$ws.Cells.Item(137, 9).Value = 990.7222  # I137: was 1081.4193
$ws.Cells.Item(137, 10).Value = 2470.1333  # J137: was 2344.2354
$ws.Cells.Item(137, 11).Value = 2972.1666  # K137: was 3244.2579
$ws.Cells.Item(137, 12).Value = 7410.3999  # L137: was 7032.706200000001
$ws.Cells.Item(137, 13).Value = -422.1666  # M137: was -694.2579000000001
$ws.Cells.Item(137, 14).Value = -12510.3999  # N137: was -12132.7062
$ws.Cells.Item(138, 8).Value = 2217.513  # H138: was 2419.8286
$ws.Cells.Item(138, 9).Value = 811.375  # I138: was 909.3
$ws.Cells.Item(138, 10).Value = 4467.3335  # J138: was 4433.8667
$ws.Cells.Item(138, 11).Value = 2434.125  # K138: was 2727.9
$ws.Cells.Item(138, 12).Value = 13402.0005  # L138: was 13301.6001
$ws.Cells.Item(138, 13).Value = 2705.875  # M138: was 2412.1
$ws.Cells.Item(138, 14).Value = -23682.0005  # N138: was -23581.6001
# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 6376.03  # H32: was 11420.107
$ws.Cells.Item(32, 9).Value = 2664.1448  # I32: was 3059.7273
$ws.Cells.Item(32, 10).Value = 18130.334  # J32: was 31856.592
$ws.Cells.Item(32, 11).Value = 2664.1448  # K32: was 3059.7273
$ws.Cells.Item(32, 12).Value = 18130.334  # L32: was 31856.592
$ws.Cells.Item(32, 13).Value = -2377.1448  # M32: was -2772.7273
$ws.Cells.Item(32, 14).Value = -18704.334  # N32: was -32430.592
$ws.Cells.Item(61, 8).Value = 758.375  # H61: was 810.75
$ws.Cells.Item(61, 9).Value = 587.5111000000001  # I61: was 652.7368
$ws.Cells.Item(61, 10).Value = 1163.0526  # J61: was 1144.3334
$ws.Cells.Item(61, 11).Value = 587.5111000000001  # K61: was 652.7368
$ws.Cells.Item(61, 12).Value = 1163.0526  # L61: was 1144.3334
$ws.Cells.Item(61, 13).Value = -375.5111000000001  # M61: was -440.7368
$ws.Cells.Item(61, 14).Value = -1587.0526  # N61: was -1568.3334
$ws.Cells.Item(74, 8).Value = 4312405.5  # H74: was 6099236.5
$ws.Cells.Item(74, 9).Value = 6946399.5  # I74: was 10870533
$ws.Cells.Item(74, 10).Value = 2234.182  # J74: was 2579.5
$ws.Cells.Item(74, 11).Value = 6946399.5  # K74: was 10870533
$ws.Cells.Item(74, 12).Value = 2234.182  # L74: was 2579.5
$ws.Cells.Item(74, 13).Value = -6945525.5  # M74: was -10869659
$ws.Cells.Item(74, 14).Value = -3982.182  # N74: was -4327.5
$ws.Cells.Item(77, 8).Value = 4312405.5  # H77: was 6099236.5
$ws.Cells.Item(77, 9).Value = 6946399.5  # I77: was 10870533
$ws.Cells.Item(77, 10).Value = 2234.182  # J77: was 2579.5
$ws.Cells.Item(77, 11).Value = 34731997.5  # K77: was 54352665
$ws.Cells.Item(77, 12).Value = 11170.91  # L77: was 12897.5
$ws.Cells.Item(77, 13).Value = -34727629.5  # M77: was -54348297
$ws.Cells.Item(77, 14).Value = -19906.91  # N77: was -21633.5
$ws.Cells.Item(132, 8).Value = 1216.5454  # H132: was 1358.3489
$ws.Cells.Item(132, 9).Value = 1068.7894  # I132: was 1138.75
$ws.Cells.Item(132, 10).Value = 2152.3333  # J132: was 2487.7144
$ws.Cells.Item(132, 11).Value = 3206.3682  # K132: was 3416.25
$ws.Cells.Item(132, 12).Value = 6456.999899999999  # L132: was 7463.1432
$ws.Cells.Item(132, 13).Value = -676.3681999999999  # M132: was -886.25
$ws.Cells.Item(132, 14).Value = -11516.9999  # N132: was -12523.1432
$ws.Cells.Item(136, 8).Value = 758.375  # H136: was 810.75
$ws.Cells.Item(136, 9).Value = 587.5111000000001  # I136: was 652.7368
$ws.Cells.Item(136, 10).Value = 1163.0526  # J136: was 1144.3334
$ws.Cells.Item(136, 11).Value = 1762.5333  # K136: was 1958.2104
$ws.Cells.Item(136, 12).Value = 3489.1578  # L136: was 3433.0002
$ws.Cells.Item(136, 13).Value = 787.4666999999999  # M136: was 591.7896000000001
$ws.Cells.Item(136, 14).Value = -8589.157800000001  # N136: was -8533.0002
# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 1638.1786  # H134: was 1898.4884
$ws.Cells.Item(134, 9).Value = 1449.5745  # I134: was 1708.2307
$ws.Cells.Item(134, 10).Value = 2623.111  # J134: was 3753.5
$ws.Cells.Item(134, 11).Value = 4348.7235  # K134: was 5124.6921
$ws.Cells.Item(134, 12).Value = 7869.333  # L134: was 11260.5
$ws.Cells.Item(134, 13).Value = -1813.7235  # M134: was -2589.6921
$ws.Cells.Item(134, 14).Value = -12939.333  # N134: was -16330.5
# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(75, 8).Value = 0  # H75: was 30000
$ws.Cells.Item(75, 10).Value = 0  # J75: was 30000
$ws.Cells.Item(75, 12).Value = 0  # L75: was 30000
$ws.Cells.Item(75, 14).ClearContents()  # N75: was -31996
$ws.Cells.Item(78, 8).Value = 0  # H78: was 30000
$ws.Cells.Item(78, 10).Value = 0  # J78: was 30000
$ws.Cells.Item(78, 12).Value = 0  # L78: was 90000
$ws.Cells.Item(78, 14).ClearContents()  # N78: was -99984
$ws.Cells.Item(132, 8).Value = 855.1818  # H132: was 1109.6744
$ws.Cells.Item(132, 9).Value = 708.8958  # I132: was 996.1818
$ws.Cells.Item(132, 10).Value = 1245.2778  # J132: was 1484.2
$ws.Cells.Item(132, 11).Value = 2126.6874  # K132: was 2988.5454
$ws.Cells.Item(132, 12).Value = 3735.8334  # L132: was 4452.6
$ws.Cells.Item(132, 13).Value = 403.3126000000002  # M132: was -458.5454
$ws.Cells.Item(132, 14).Value = -8795.8334  # N132: was -9512.6
$ws.Cells.Item(134, 8).Value = 2544.182  # H134: was 868.0755
$ws.Cells.Item(134, 9).Value = 2902.9614  # I134: was 702.7838
$ws.Cells.Item(134, 10).Value = 1211.5714  # J134: was 1250.3125
$ws.Cells.Item(134, 11).Value = 8708.8842  # K134: was 2108.3514
$ws.Cells.Item(134, 12).Value = 3634.7142  # L134: was 3750.9375
$ws.Cells.Item(134, 13).Value = -6173.8842  # M134: was 426.6486
$ws.Cells.Item(134, 14).Value = -8704.7142  # N134: was -8820.9375
# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 458.0263  # H5: was 623.1875
$ws.Cells.Item(5, 9).Value = 419.58334  # I5: was 533.4643
$ws.Cells.Item(5, 10).Value = 1150  # J5: was 1251.25
$ws.Cells.Item(5, 11).Value = 1258.75002  # K5: was 1600.3929
$ws.Cells.Item(5, 12).Value = 3450  # L5: was 3753.75
$ws.Cells.Item(5, 13).Value = -1146.75002  # M5: was -1488.3929
$ws.Cells.Item(5, 14).Value = -3674  # N5: was -3977.75
$ws.Cells.Item(12, 8).Value = 40.35  # H12: was 57.75
$ws.Cells.Item(12, 9).Value = 8.25  # I12: was 11
$ws.Cells.Item(12, 10).Value = 48.375  # J12: was 62
$ws.Cells.Item(12, 11).Value = 24.75  # K12: was 33
$ws.Cells.Item(12, 12).Value = 145.125  # L12: was 186
$ws.Cells.Item(12, 13).Value = 148.25  # M12: was 140
$ws.Cells.Item(12, 14).Value = -491.125  # N12: was -532
$ws.Cells.Item(14, 8).Value = 55.545456  # H14: was 68.666664
$ws.Cells.Item(14, 9).Value = 55.545456  # I14: was 68.666664
$ws.Cells.Item(14, 11).Value = 166.636368  # K14: was 205.999992
$ws.Cells.Item(14, 13).Value = 6.363631999999996  # M14: was -32.99999199999999
$ws.Cells.Item(33, 8).Value = 8587.166999999999  # H33: was 4941.091
$ws.Cells.Item(33, 9).Value = 384.2  # I33: was 634.0833
$ws.Cells.Item(33, 10).Value = 14446.429  # J33: was 10109.5
$ws.Cells.Item(33, 11).Value = 2305.2  # K33: was 3804.4998
$ws.Cells.Item(33, 12).Value = 86678.57399999999  # L33: was 60657
$ws.Cells.Item(33, 13).Value = -2022.2  # M33: was -3521.4998
$ws.Cells.Item(33, 14).Value = -87244.57399999999  # N33: was -61223
$ws.Cells.Item(38, 8).Value = 216.78572  # H38: was 296.875
$ws.Cells.Item(38, 9).Value = 176.75  # I38: was 298.5
$ws.Cells.Item(38, 10).Value = 270.16666  # J38: was 295.25
$ws.Cells.Item(38, 11).Value = 530.25  # K38: was 895.5
$ws.Cells.Item(38, 12).Value = 810.4999799999999  # L38: was 885.75
$ws.Cells.Item(38, 13).Value = -183.25  # M38: was -548.5
$ws.Cells.Item(38, 14).Value = -1504.49998  # N38: was -1579.75
$ws.Cells.Item(40, 8).Value = 7732.923  # H40: was 5923.5884
$ws.Cells.Item(40, 9).Value = 129.8  # I40: was 41.125
$ws.Cells.Item(40, 10).Value = 12484.875  # J40: was 11152.444
$ws.Cells.Item(40, 11).Value = 519.2  # K40: was 164.5
$ws.Cells.Item(40, 12).Value = 49939.5  # L40: was 44609.776
$ws.Cells.Item(40, 13).Value = -450.2  # M40: was -95.5
$ws.Cells.Item(40, 14).Value = -50077.5  # N40: was -44747.776
$ws.Cells.Item(68, 8).Value = 610.2222  # H68: was 807.05554
$ws.Cells.Item(68, 9).Value = 330  # I68: was 400
$ws.Cells.Item(68, 10).Value = 666.26666  # J68: was 831
$ws.Cells.Item(68, 11).Value = 990  # K68: was 1200
$ws.Cells.Item(68, 12).Value = 1998.79998  # L68: was 2493
$ws.Cells.Item(68, 13).Value = -179  # M68: was -389
$ws.Cells.Item(68, 14).Value = -3620.79998  # N68: was -4115
$ws.Cells.Item(71, 8).Value = 610.2222  # H71: was 807.05554
$ws.Cells.Item(71, 9).Value = 330  # I71: was 400
$ws.Cells.Item(71, 10).Value = 666.26666  # J71: was 831
$ws.Cells.Item(71, 11).Value = 2970  # K71: was 3600
$ws.Cells.Item(71, 12).Value = 5996.39994  # L71: was 7479
$ws.Cells.Item(71, 13).Value = 1086  # M71: was 456
$ws.Cells.Item(71, 14).Value = -14108.39994  # N71: was -15591
$ws.Cells.Item(80, 8).Value = 1676.6666  # H80: was 2000
$ws.Cells.Item(80, 9).Value = 1000  # I80: was 0
$ws.Cells.Item(80, 10).Value = 1812  # J80: was 2000
$ws.Cells.Item(80, 11).Value = 3000  # K80: was 0
$ws.Cells.Item(80, 12).Value = 5436  # L80: was 6000
$ws.Cells.Item(80, 13).Value = -2064  # M80: was None
$ws.Cells.Item(80, 14).Value = -7308  # N80: was -7872
$ws.Cells.Item(83, 8).Value = 1676.6666  # H83: was 2000
$ws.Cells.Item(83, 9).Value = 1000  # I83: was 0
$ws.Cells.Item(83, 10).Value = 1812  # J83: was 2000
$ws.Cells.Item(83, 11).Value = 9000  # K83: was 0
$ws.Cells.Item(83, 12).Value = 16308  # L83: was 18000
$ws.Cells.Item(83, 13).Value = -4320  # M83: was None
$ws.Cells.Item(83, 14).Value = -25668  # N83: was -27360
$ws.Cells.Item(86, 8).Value = 820.9286  # H86: was 201.73334
$ws.Cells.Item(86, 9).Value = 812.5  # I86: was 203
$ws.Cells.Item(86, 10).Value = 832.1667  # J86: was 196.66667
$ws.Cells.Item(86, 11).Value = 2437.5  # K86: was 609
$ws.Cells.Item(86, 12).Value = 2496.5001  # L86: was 590.00001
$ws.Cells.Item(86, 13).Value = -1251.5  # M86: was 577
$ws.Cells.Item(86, 14).Value = -4868.5001  # N86: was -2962.00001
$ws.Cells.Item(89, 8).Value = 820.9286  # H89: was 201.73334
$ws.Cells.Item(89, 9).Value = 812.5  # I89: was 203
$ws.Cells.Item(89, 10).Value = 832.1667  # J89: was 196.66667
$ws.Cells.Item(89, 11).Value = 7312.5  # K89: was 1827
$ws.Cells.Item(89, 12).Value = 7489.5003  # L89: was 1770.00003
$ws.Cells.Item(89, 13).Value = -1384.5  # M89: was 4101
$ws.Cells.Item(89, 14).Value = -19345.5003  # N89: was -13626.00003
$ws.Cells.Item(97, 8).Value = 348.57144  # H97: was 495.2381
$ws.Cells.Item(97, 9).Value = 389.9  # I97: was 412.9
$ws.Cells.Item(97, 10).Value = 311  # J97: was 570.0909
$ws.Cells.Item(97, 11).Value = 1169.7  # K97: was 1238.7
$ws.Cells.Item(97, 12).Value = 933  # L97: was 1710.2727
$ws.Cells.Item(97, 13).Value = -673.6999999999998  # M97: was -742.6999999999998
$ws.Cells.Item(97, 14).Value = -1925  # N97: was -2702.2727
$ws.Cells.Item(98, 8).Value = 397.85  # H98: was 408.10526
$ws.Cells.Item(98, 9).Value = 410.5  # I98: was 433.55554
$ws.Cells.Item(98, 11).Value = 1231.5  # K98: was 1300.66662
$ws.Cells.Item(98, 13).Value = 266.5  # M98: was 197.33338
$ws.Cells.Item(107, 8).Value = 319.875  # H107: was 223.78572
$ws.Cells.Item(107, 9).Value = 296.57144  # I107: was 194.3077
$ws.Cells.Item(107, 10).Value = 338  # J107: was 249.33333
$ws.Cells.Item(107, 11).Value = 889.71432  # K107: was 582.9231
$ws.Cells.Item(107, 12).Value = 1014  # L107: was 747.99999
$ws.Cells.Item(107, 13).Value = 1030.28568  # M107: was 1337.0769
$ws.Cells.Item(107, 14).Value = -4854  # N107: was -4587.99999
$ws.Cells.Item(113, 8).Value = 1277.7858  # H113: was 1872.75
$ws.Cells.Item(113, 9).Value = 544.3333  # I113: was 1000
$ws.Cells.Item(113, 10).Value = 1827.875  # J113: was 1997.4286
$ws.Cells.Item(113, 11).Value = 1632.9999  # K113: was 3000
$ws.Cells.Item(113, 12).Value = 5483.625  # L113: was 5992.2858
$ws.Cells.Item(113, 13).Value = 537.0001  # M113: was -830
$ws.Cells.Item(113, 14).Value = -9823.625  # N113: was -10332.2858
$ws.Cells.Item(131, 8).Value = 1516057.1  # H131: was 1235506.1
$ws.Cells.Item(131, 9).Value = 4762219  # I131: was 5128528
$ws.Cells.Item(131, 10).Value = 1181.5333  # J131: was 1133.317
$ws.Cells.Item(131, 11).Value = 14286657  # K131: was 15385584
$ws.Cells.Item(131, 12).Value = 3544.5999  # L131: was 3399.951
$ws.Cells.Item(131, 13).Value = -14281617  # M131: was -15380544
$ws.Cells.Item(131, 14).Value = -13624.5999  # N131: was -13479.951
$ws.Cells.Item(132, 8).Value = 748  # H132: was 712.5
$ws.Cells.Item(132, 9).Value = 558.2222  # I132: was 733.3333
$ws.Cells.Item(132, 10).Value = 1175  # J132: was 700
$ws.Cells.Item(132, 11).Value = 5023.999800000001  # K132: was 6599.9997
$ws.Cells.Item(132, 12).Value = 10575  # L132: was 6300
$ws.Cells.Item(132, 13).Value = -2493.999800000001  # M132: was -4069.9997
$ws.Cells.Item(132, 14).Value = -15635  # N132: was -11360
$ws.Cells.Item(135, 8).Value = 458.0263  # H135: was 623.1875
$ws.Cells.Item(135, 9).Value = 419.58334  # I135: was 533.4643
$ws.Cells.Item(135, 10).Value = 1150  # J135: was 1251.25
$ws.Cells.Item(135, 11).Value = 3776.25006  # K135: was 4801.178699999999
$ws.Cells.Item(135, 12).Value = 10350  # L135: was 11261.25
$ws.Cells.Item(135, 13).Value = -1241.25006  # M135: was -2266.178699999999
$ws.Cells.Item(135, 14).Value = -15420  # N135: was -16331.25
# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 100.05882  # H2: was 136.61539
$ws.Cells.Item(2, 9).Value = 69.666664  # I2: was 95.40000000000001
$ws.Cells.Item(2, 10).Value = 173  # J2: was 274
$ws.Cells.Item(2, 11).Value = 69.666664  # K2: was 95.40000000000001
$ws.Cells.Item(2, 12).Value = 173  # L2: was 274
$ws.Cells.Item(2, 13).Value = 43.333336  # M2: was 17.59999999999999
$ws.Cells.Item(2, 14).Value = -399  # N2: was -500
# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 7271184  # H132: was 9620011
$ws.Cells.Item(132, 9).Value = 11165287  # I132: was 16031899
$ws.Cells.Item(132, 10).Value = 2191.3667  # J132: was 2180.3462
$ws.Cells.Item(132, 11).Value = 33495861  # K132: was 48095697
$ws.Cells.Item(132, 12).Value = 6574.1001  # L132: was 6541.0386
$ws.Cells.Item(132, 13).Value = -33493331  # M132: was -48093167
$ws.Cells.Item(132, 14).Value = -11634.1001  # N132: was -11601.0386
# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 1097.6034  # H136: was 841.9394
$ws.Cells.Item(136, 9).Value = 725.63416  # I136: was 440.89584
$ws.Cells.Item(136, 10).Value = 1994.7059  # J136: was 1911.3889
$ws.Cells.Item(136, 11).Value = 2176.90248  # K136: was 1322.68752
$ws.Cells.Item(136, 12).Value = 5984.1177  # L136: was 5734.1667
$ws.Cells.Item(136, 13).Value = 373.0975200000003  # M136: was 1227.31248
$ws.Cells.Item(136, 14).Value = -11084.1177  # N136: was -10834.1667
